# Applies the benchmark-stats correction described in the commit:
# "Fixed README.md stats and docx preparation for all DaCapo - JDK 21 - Shenandoah GC tests"
#
# The document is a single-column table. We:
#   1. Rewrite a handful of single-value cells (rows 1,2,3,4,6,7,10,11,12).
#   2. Delete the row that held the now-orphaned "0.00003" value (row 8).
#   3. Insert a new row (after the row that used to hold "0.00479", now
#      carrying "0.00043") holding the new "0.04037" value.
#   4. Collapse three tab-separated multi-value rows (44,45,46) down to the
#      single summary value that used to live at the top of the table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $rowIndex, $text) {
    $cell = $table.Rows.Item($rowIndex).Cells.Item(1)
    $cell.Range.Text = $text
}

# --- Simple value replacements (row indices unaffected so far) ---
Set-CellText $t 1 "0M"
Set-CellText $t 2 "0M"
Set-CellText $t 3 "0M"
Set-CellText $t 4 "212"
Set-CellText $t 6 "0.00055"
Set-CellText $t 7 "0.00019"

# --- Remove the row that used to hold "0.00003" (row 8) ---
$t.Rows.Item(8).Delete()

# After the delete, the old rows 10/11/12 are now rows 9/10/11.
Set-CellText $t 9  "0.00026"
Set-CellText $t 10 "0.00034"
Set-CellText $t 11 "0.00043"

# --- Insert a new row right after that one (now row 11), with the new value ---
$newRow = $t.Rows.Add($t.Rows.Item(12))
$newCell = $newRow.Cells.Item(1)
$newCell.Range.Text = "0.04037"

# The delete (-1) + insert (+1) before the tail rows cancel out, so rows
# 44/45/46 are still at 44/45/46.
Set-CellText $t 44 "99.98"
Set-CellText $t 45 "0.04"
Set-CellText $t 46 "229"

Write-Output ("Final row count: " + $t.Rows.Count)
